$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Rename the new issue sheets (new exports from the "MySql for Excel" tool) ---
$ws2.Name = "2151104"
$ws3.Name = "20151105"

# --- Populate "2151104" (Deploy ally makes game stuck / Bugfix) ---
$ws2.Range("A1").Value = "Issue Summary"
$ws2.Range("B1").Value = "Issue Type"
$ws2.Range("C1").Value = "Resolved?"
$ws2.Range("D1").Value = "Detail"
$ws2.Range("E1").Value = "Comment"

$ws2.Range("A2").Value = "Deploy ally makes game stuck"
$ws2.Range("B2").Value = "Bugfix"
$ws2.Range("D2").Value = "When there is not enough allies, the bot clicks the next page but it gets to the last page and still tries to find the next page button, but I can't. so it gets stuck"

# Header formatting + column sizing, matching the exported-table look of the other sheets
$ws2.Range("A1:E1").HorizontalAlignment = -4108
$ws2.Columns("A").ColumnWidth = 27
$ws2.Columns("B").ColumnWidth = 9.5
$ws2.Columns("D").ColumnWidth = 141.83333333333334

$ws2.Range("A1:XFD1").Select() | Out-Null

# --- Populate "20151105" (Add script for Frontlines Event / Feature Request) ---
$ws3.Range("A1").Value = "Issue Summary"
$ws3.Range("B1").Value = "Issue Type"
$ws3.Range("C1").Value = "Resolved?"
$ws3.Range("D1").Value = "Detail"
$ws3.Range("E1").Value = "Comment"

$ws3.Range("A2").Value = "Add script for Frontlines Event"
$ws3.Range("B2").Value = "Feature Request"
$ws3.Range("C2").Value = "Yes"
$ws3.Range("D2").Value = "Want to automated the frontlines event."

# Header formatting + column sizing
$ws3.Range("A1:E1").HorizontalAlignment = -4108
$ws3.Columns("A").ColumnWidth = 27
$ws3.Columns("B").ColumnWidth = 9.5
$ws3.Columns("D").ColumnWidth = 141.83333333333334

$ws3.Range("D2").Select() | Out-Null

# --- Sheet1: clear the old "working" selection/topLeftCell, select the header row ---
$ws1.Range("A1:XFD1").Select() | Out-Null

# --- Make "20151105" the active tab ---
$ws3.Activate() | Out-Null
